$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serping1"
$ws.Range("C2").Value = "Lrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 6.159891999999999
$ws.Range("H2").Value = 18.479676
$ws.Range("I2").Value = 0.007079533182016282
$ws.Range("J2").Value = 0.007079533182016282
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 145.5939609605186
$ws.Range("R2").Value = 1310.345648644668
$ws.Range("S2").Value = 0.0004833795231971438
$ws.Range("T2").Value = 0.0004833795231971438

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serping1"
$ws.Range("C3").Value = "Lrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 6.159891999999999
$ws.Range("H3").Value = 18.479676
$ws.Range("I3").Value = 0.007079533182016282
$ws.Range("J3").Value = 0.007079533182016282
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 1116.716911467069
$ws.Range("R3").Value = 10050.45220320362
$ws.Range("S3").Value = 0.003707558230093888
$ws.Range("T3").Value = 0.003707558230093887

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Serping1"
$ws.Range("C4").Value = "Lrp1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 6.159891999999999
$ws.Range("H4").Value = 18.479676
$ws.Range("I4").Value = 0.007079533182016282
$ws.Range("J4").Value = 0.007079533182016282
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 684.3673665543292
$ws.Range("R4").Value = 6159.306298988962
$ws.Range("S4").Value = 0.002272135253099019
$ws.Range("T4").Value = 0.002272135253099019

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Serping1"
$ws.Range("C5").Value = "Lrp1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 6.159891999999999
$ws.Range("H5").Value = 18.479676
$ws.Range("I5").Value = 0.007079533182016282
$ws.Range("J5").Value = 0.007079533182016282
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 185.6778668450853
$ws.Range("R5").Value = 1671.100801605768
$ws.Range("S5").Value = 0.0006164601756262321
$ws.Range("T5").Value = 0.0006164601756262321

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Serping1"
$ws.Range("C6").Value = "Lrp1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 771.4717003333334
$ws.Range("H6").Value = 2314.415101
$ws.Range("I6").Value = 0.8866485810946614
$ws.Range("J6").Value = 0.8866485810946614
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 18234.34901463797
$ws.Range("R6").Value = 164109.1411317417
$ws.Range("S6").Value = 0.06053898715549178
$ws.Range("T6").Value = 0.06053898715549178

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Serping1"
$ws.Range("C7").Value = "Lrp1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 771.4717003333334
$ws.Range("H7").Value = 2314.415101
$ws.Range("I7").Value = 0.8866485810946614
$ws.Range("J7").Value = 0.8866485810946614
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 139858.8635126214
$ws.Range("R7").Value = 1258729.771613593
$ws.Range("S7").Value = 0.4643387013693383
$ws.Range("T7").Value = 0.4643387013693382

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Serping1"
$ws.Range("C8").Value = "Lrp1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 771.4717003333334
$ws.Range("H8").Value = 2314.415101
$ws.Range("I8").Value = 0.8866485810946614
$ws.Range("J8").Value = 0.8866485810946614
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 85710.92738773895
$ws.Range("R8").Value = 771398.3464896504
$ws.Range("S8").Value = 0.2845647370271441
$ws.Range("T8").Value = 0.2845647370271441

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Serping1"
$ws.Range("C9").Value = "Lrp1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 771.4717003333334
$ws.Range("H9").Value = 2314.415101
$ws.Range("I9").Value = 0.8866485810946614
$ws.Range("J9").Value = 0.8866485810946614
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 23254.50180770121
$ws.Range("R9").Value = 209290.5162693109
$ws.Range("S9").Value = 0.07720615554268721
$ws.Range("T9").Value = 0.07720615554268721

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Serping1"
$ws.Range("C10").Value = "Lrp1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4146736666666667
$ws.Range("H10").Value = 1.244021
$ws.Range("I10").Value = 0.0004765823788590816
$ws.Range("J10").Value = 0.0004765823788590817
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 9.80114288302811
$ws.Range("R10").Value = 88.210285947253
$ws.Range("S10").Value = 0.00003254030416048605
$ws.Range("T10").Value = 0.00003254030416048605

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Serping1"
$ws.Range("C11").Value = "Lrp1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.4146736666666667
$ws.Range("H11").Value = 1.244021
$ws.Range("I11").Value = 0.0004765823788590816
$ws.Range("J11").Value = 0.0004765823788590817
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 75.1755219582949
$ws.Range("R11").Value = 676.579697624654
$ws.Range("S11").Value = 0.0002495866430212104
$ws.Range("T11").Value = 0.0002495866430212104

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Serping1"
$ws.Range("C12").Value = "Lrp1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4146736666666667
$ws.Range("H12").Value = 1.244021
$ws.Range("I12").Value = 0.0004765823788590816
$ws.Range("J12").Value = 0.0004765823788590817
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 46.07047091671322
$ws.Range("R12").Value = 414.6342382504189
$ws.Range("S12").Value = 0.0001529563597162361
$ws.Range("T12").Value = 0.0001529563597162361

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Serping1"
$ws.Range("C13").Value = "Lrp1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4146736666666667
$ws.Range("H13").Value = 1.244021
$ws.Range("I13").Value = 0.0004765823788590816
$ws.Range("J13").Value = 0.0004765823788590817
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 12.49952464483089
$ws.Range("R13").Value = 112.495721803478
$ws.Range("S13").Value = 0.00004149907196114915
$ws.Range("T13").Value = 0.00004149907196114915

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Serping1"
$ws.Range("C14").Value = "Lrp1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 92.05234666666666
$ws.Range("H14").Value = 276.15704
$ws.Range("I14").Value = 0.1057953033444633
$ws.Range("J14").Value = 0.1057953033444633
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 2175.730640555191
$ws.Range("R14").Value = 19581.57576499672
$ws.Range("S14").Value = 0.007223538893362341
$ws.Range("T14").Value = 0.007223538893362341

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Serping1"
$ws.Range("C15").Value = "Lrp1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 92.05234666666666
$ws.Range("H15").Value = 276.15704
$ws.Range("I15").Value = 0.1057953033444633
$ws.Range("J15").Value = 0.1057953033444633
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 16688.02184565833
$ws.Range("R15").Value = 150192.196610925
$ws.Range("S15").Value = 0.05540510052505072
$ws.Range("T15").Value = 0.0554051005250507

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Serping1"
$ws.Range("C16").Value = "Lrp1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 92.05234666666666
$ws.Range("H16").Value = 276.15704
$ws.Range("I16").Value = 0.1057953033444633
$ws.Range("J16").Value = 0.1057953033444633
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 10227.06600593206
$ws.Range("R16").Value = 92043.59405338854
$ws.Range("S16").Value = 0.03395439108215294
$ws.Range("T16").Value = 0.03395439108215294

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Serping1"
$ws.Range("C17").Value = "Lrp1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 92.05234666666666
$ws.Range("H17").Value = 276.15704
$ws.Range("I17").Value = 0.1057953033444633
$ws.Range("J17").Value = 0.1057953033444633
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 2774.737506298969
$ws.Range("R17").Value = 24972.63755669072
$ws.Range("S17").Value = 0.009212272843897284
$ws.Range("T17").Value = 0.009212272843897284

